# Update the cryptocurrency price table to reflect the latest scrape.
# A new coin (WrappedliquidstakedEther2.0) is inserted at row 17, shifting
# the remaining rows down by one (the final row, RenderToken, drops off the
# bottom of the existing 50-row table). All Price (D) and Volume(1h) (E)
# values are refreshed with updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.427.41'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '1.849.61'
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '240.53'
$c.ClearFormats()
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  +0.05%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07704'
$c.ClearFormats()
$ws.Range("E8").Value = '  +2.36%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2917'
$c.ClearFormats()
$ws.Range("E9").Value = '  -0.31%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '24.83'
$c.ClearFormats()
$ws.Range("E10").Value = '  +1.23%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07753'
$c.ClearFormats()
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").Value = '1.852.37'
$ws.Range("E12").Value = '  +0.51%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.029'
$c.ClearFormats()
$ws.Range("E13").Value = '  +0.55%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.00001079'
$c.ClearFormats()
$ws.Range("E14").Value = '  +3.58%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.6815'
$c.ClearFormats()
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '2.121.26'
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '6.178'
$c.ClearFormats()
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '29.469.43'
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '228.76'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '12.38'
$c.ClearFormats()
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '7.438'
$c.ClearFormats()
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '157.52'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.1377'
$c.ClearFormats()
$ws.Range("E26").Value = '  -1.19%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.420'
$c.ClearFormats()
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '17.72'
$c.ClearFormats()
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.346'
$c.ClearFormats()
$ws.Range("E29").Value = '  +4.99%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.462'
$c.ClearFormats()
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.05646'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.124'
$c.ClearFormats()
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.038'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.843'
$c.ClearFormats()
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.163'
$c.ClearFormats()
$ws.Range("E35").Value = '  +0.58%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.7025'
$c.ClearFormats()
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.594'
$c.ClearFormats()
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.226.69'
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.765'
$c.ClearFormats()
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.01788'
$c.ClearFormats()
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.444'
$c.ClearFormats()
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.9066'
$c.ClearFormats()
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '2.027.05'
$ws.Range("E44").Value = '  +1.38%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '101.85'
$c.ClearFormats()
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '65.97'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000120'
$c.ClearFormats()
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.171'
$c.ClearFormats()
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.4018'
$c.ClearFormats()
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.1160'
$c.ClearFormats()
$ws.Range("E50").Value = '  +3.42%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '9.018'
$c.ClearFormats()
$ws.Range("E51").Value = '  +0.62%  '
